$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CreatedAt timestamp in A1
$ws.Range("A1").Value = "CreatedAt: 2025-08-16T18:07:43"

# Update data values that changed in this report refresh
$ws.Range("V4").Value = 117.82
$ws.Range("W4").Value = 105
$ws.Range("X4").Value = 54.41
$ws.Range("Y4").Value = 68.45999999999999
$ws.Range("Z4").Value = 70.45999999999999
$ws.Range("V5").Value = -25.63
$ws.Range("W5").Value = -34.43
$ws.Range("V6").Value = -11.05
$ws.Range("W6").Value = -11.71
$ws.Range("X6").Value = -3.54
$ws.Range("Y6").Value = -3.42
$ws.Range("Z6").Value = -2.82
$ws.Range("V9").Value = 113.18
$ws.Range("W9").Value = 99.8
$ws.Range("X9").Value = 54.92
$ws.Range("Y9").Value = 70.33
$ws.Range("Z9").Value = 72.77
$ws.Range("V10").Value = -25.63
$ws.Range("W10").Value = -34.43
$ws.Range("V11").Value = -15.69
$ws.Range("W11").Value = -16.91
$ws.Range("X11").Value = -3.02
$ws.Range("Y11").Value = -1.55
$ws.Range("Z11").Value = -0.51
$ws.Range("V14").Value = 122.45
$ws.Range("W14").Value = 134.23
$ws.Range("X14").Value = 54.92
$ws.Range("Y14").Value = 70.33
$ws.Range("Z14").Value = 72.84
$ws.Range("V15").Value = -16.36
$ws.Range("V16").Value = -15.69
$ws.Range("W16").Value = -16.91
$ws.Range("X16").Value = -3.02
$ws.Range("Y16").Value = -1.55
$ws.Range("Z16").Value = -0.44
$ws.Range("W19").Value = 137.9
$ws.Range("X19").Value = 54.05
$ws.Range("Y19").Value = 68.06999999999999
$ws.Range("Z19").Value = 70.12
$ws.Range("V21").Value = -12.1
$ws.Range("W21").Value = -13.24
$ws.Range("X21").Value = -3.89
$ws.Range("Y21").Value = -3.81
$ws.Range("Z21").Value = -3.16
$ws.Range("V24").Value = 116.77
$ws.Range("W24").Value = 103.47
$ws.Range("X24").Value = 54.05
$ws.Range("Y24").Value = 68.06999999999999
$ws.Range("Z24").Value = 70.12
$ws.Range("V25").Value = -25.63
$ws.Range("W25").Value = -34.43
$ws.Range("V26").Value = -12.1
$ws.Range("W26").Value = -13.24
$ws.Range("X26").Value = -3.89
$ws.Range("Y26").Value = -3.81
$ws.Range("Z26").Value = -3.16
$ws.Range("W29").Value = 136.41
$ws.Range("X29").Value = 53.6
$ws.Range("Y29").Value = 67.62
$ws.Range("Z29").Value = 69.52
$ws.Range("V31").Value = -12.89
$ws.Range("W31").Value = -14.73
$ws.Range("X31").Value = -4.34
$ws.Range("Y31").Value = -4.26
$ws.Range("Z31").Value = -3.75
$ws.Range("V34").Value = 120
$ws.Range("W34").Value = 131.2
$ws.Range("X34").Value = 55.82
$ws.Range("Y34").Value = 71.81
$ws.Range("Z34").Value = 74.62
$ws.Range("V35").Value = -16.36
$ws.Range("V36").Value = -18.14
$ws.Range("W36").Value = -19.94
$ws.Range("X36").Value = -2.12
$ws.Range("Y36").Value = -0.07000000000000001
$ws.Range("Z36").Value = 1.34
$ws.Range("V39").Value = 117.82
$ws.Range("W39").Value = 105
$ws.Range("X39").Value = 54.41
$ws.Range("Y39").Value = 68.45999999999999
$ws.Range("Z39").Value = 70.45999999999999
$ws.Range("V40").Value = -25.63
$ws.Range("W40").Value = -34.43
$ws.Range("V41").Value = -11.05
$ws.Range("W41").Value = -11.71
$ws.Range("X41").Value = -3.54
$ws.Range("Y41").Value = -3.42
$ws.Range("Z41").Value = -2.82
$ws.Range("V44").Value = 150.58
$ws.Range("W44").Value = 100.62
$ws.Range("X44").Value = 54.59
$ws.Range("Y44").Value = 69.53
$ws.Range("Z44").Value = 71.84999999999999
$ws.Range("W45").Value = -34.43
$ws.Range("V46").Value = -3.92
$ws.Range("W46").Value = -16.1
$ws.Range("X46").Value = -3.35
$ws.Range("Y46").Value = -2.35
$ws.Range("Z46").Value = -1.43
$ws.Range("V49").Value = 153.42
$ws.Range("W49").Value = 157.93
$ws.Range("X49").Value = 60.11
$ws.Range("Y49").Value = 74.95
$ws.Range("Z49").Value = 75.78
$ws.Range("V51").Value = -1.07
$ws.Range("W51").Value = 6.79
$ws.Range("X51").Value = 2.16
$ws.Range("Y51").Value = 3.07
$ws.Range("Z51").Value = 2.5
$ws.Range("V54").Value = 143.72
$ws.Range("W54").Value = 141.67
$ws.Range("X54").Value = 55.93
$ws.Range("Y54").Value = 70.54000000000001
$ws.Range("Z54").Value = 73.5
$ws.Range("V56").Value = -10.51
$ws.Range("W56").Value = -9.359999999999999
$ws.Range("X56").Value = -2.01
$ws.Range("Y56").Value = -1.34
$ws.Range("Z56").Value = 0.22
$ws.Range("V57").Value = -0.26
$ws.Range("W57").Value = -0.11
$ws.Range("V59").Value = 155.12
$ws.Range("W59").Value = 157.77
$ws.Range("X59").Value = 59.98
$ws.Range("Y59").Value = 74.56
$ws.Range("Z59").Value = 75.45999999999999
$ws.Range("V61").Value = 0.62
$ws.Range("W61").Value = 6.63
$ws.Range("X61").Value = 2.04
$ws.Range("Y61").Value = 2.68
$ws.Range("Z61").Value = 2.19
$ws.Range("V64").Value = 157.97
$ws.Range("W64").Value = 159.27
$ws.Range("X64").Value = 60.8
$ws.Range("Y64").Value = 75.5
$ws.Range("Z64").Value = 76.25
$ws.Range("V66").Value = 3.48
$ws.Range("W66").Value = 8.119999999999999
$ws.Range("X66").Value = 2.86
$ws.Range("Y66").Value = 3.62
$ws.Range("Z66").Value = 2.97
$ws.Range("V69").Value = 155.43
$ws.Range("W69").Value = 160.11
$ws.Range("X69").Value = 60.61
$ws.Range("Y69").Value = 75.43000000000001
$ws.Range("Z69").Value = 76.56999999999999
$ws.Range("V71").Value = 0.93
$ws.Range("W71").Value = 8.970000000000001
$ws.Range("X71").Value = 2.67
$ws.Range("Y71").Value = 3.55
$ws.Range("Z71").Value = 3.29
$ws.Range("V74").Value = 155.59
$ws.Range("W74").Value = 156.95
$ws.Range("X74").Value = 59.8
$ws.Range("Y74").Value = 74.26000000000001
$ws.Range("Z74").Value = 75.15000000000001
$ws.Range("V76").Value = 1.09
$ws.Range("W76").Value = 5.81
$ws.Range("X76").Value = 1.85
$ws.Range("Z76").Value = 1.88
$ws.Range("V79").Value = 156.28
$ws.Range("W79").Value = 157.51
$ws.Range("X79").Value = 60.13
$ws.Range("Y79").Value = 74.68000000000001
$ws.Range("Z79").Value = 75.62
$ws.Range("V81").Value = 1.78
$ws.Range("W81").Value = 6.36
$ws.Range("X81").Value = 2.19
$ws.Range("Y81").Value = 2.8
$ws.Range("Z81").Value = 2.35
$ws.Range("W84").Value = 138.16
$ws.Range("X84").Value = 55.98
$ws.Range("Y84").Value = 70.2
$ws.Range("Z84").Value = 72.48
$ws.Range("V86").Value = -13.27
$ws.Range("W86").Value = -12.99
$ws.Range("X86").Value = -1.96
$ws.Range("Y86").Value = -1.68
$ws.Range("Z86").Value = -0.8
$ws.Range("W89").Value = 136.41
$ws.Range("X89").Value = 53.6
$ws.Range("Y89").Value = 67.62
$ws.Range("Z89").Value = 69.52
$ws.Range("V91").Value = -12.89
$ws.Range("W91").Value = -14.73
$ws.Range("X91").Value = -4.34
$ws.Range("Y91").Value = -4.26
$ws.Range("Z91").Value = -3.75
